# Insert a new data row at row 136 (shifts existing rows 136:264 down to 137:265)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with its data
$ws.Range("A136").Value = 10
$ws.Range("B136").Value = "Vega Modelo de Temuco"
$ws.Range("C136").Value = "La Araucanía"
$ws.Range("D136").Value = 44589
$ws.Range("E136").Value = 9
$ws.Range("F136").Value = 100112009
$ws.Range("G136").Value = "Acelga"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 85
$ws.Range("K136").Value = 7000
$ws.Range("L136").Value = 8000
$ws.Range("M136").Value = 7529
$ws.Range("N136").Value = "$/docena de atados (12 kilos)"
$ws.Range("O136").Value = "Provincia de Cautín"
$ws.Range("P136").Value = 627
$ws.Range("Q136").Value = 12
$ws.Range("R136").Value = "Hortaliza"
